$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 419, shifting existing rows 419:445 down to 420:446,
# and populate it with this week's new price report for Brócoli.
$ws.Rows.Item(419).EntireRow.Insert()

$ws.Range("A419").Value = 10
$ws.Range("B419").Value = "Vega Modelo de Temuco"
$ws.Range("C419").Value = "La Araucanía"
$ws.Range("D419").Value = 44746
$ws.Range("E419").Value = 9
$ws.Range("F419").Value = 100112023
$ws.Range("G419").Value = "Brócoli"
$ws.Range("H419").Value = "Sin especificar"
$ws.Range("I419").Value = "Primera"
$ws.Range("J419").Value = 2150
$ws.Range("K419").Value = 1000
$ws.Range("L419").Value = 1200
$ws.Range("M419").Value = 1084
$ws.Range("N419").Value = "$/unidad"
$ws.Range("O419").Value = "Región del Maule"
$ws.Range("P419").Value = 1084
$ws.Range("Q419").Value = 1
$ws.Range("R419").Value = "Hortaliza"
